# Insert a new data row at row 96 (pushes existing rows 96-118 down to 97-119)
# and populate it with a new weekly price record for "Poroto verde".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(96).Insert()

$ws.Range("A96").Value = 5
$ws.Range("B96").Value = "Macroferia Regional de Talca"
$ws.Range("C96").Value = "Maule"
$ws.Range("D96").Value = 44543
$ws.Range("E96").Value = 7
$ws.Range("F96").Value = 100112031
$ws.Range("G96").Value = "Poroto verde"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 400
$ws.Range("K96").Value = 12000
$ws.Range("L96").Value = 12000
$ws.Range("M96").Value = 12000
$ws.Range("N96").Value = "`$/saco 25 kilos"
$ws.Range("O96").Value = "Región del Maule"
$ws.Range("P96").Value = 480
$ws.Range("Q96").Value = 25
$ws.Range("R96").Value = "Hortaliza"
